$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = "Amira Al Jabri"
$wsSummary.Range("B4").Value = 2281.69
$wsSummary.Range("B6").Value = 113019
$wsSummary.Range("B7").Value = 46463
$wsSummary.Range("B8").Value = 66556
$wsSummary.Range("B9").Value = 2.43

# --- Assets sheet ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Range("C2").Value = 111120
$wsAssets.Range("C3").Value = 1899
$wsAssets.Range("C4").Value = 113019

# --- Liabilities sheet ---
$wsLiabilities = $wb.Worksheets.Item("Liabilities")
$wsLiabilities.Range("C2").Value = 46463
$wsLiabilities.Range("D2").Value = 2323
$wsLiabilities.Range("C3").Value = 46463
